$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.473.83"
$ws.Range("E2").Value = "  +2.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.646.97"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.02"
$ws.Range("E5").Value = "  +1.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.42"
$ws.Range("E6").Value = "  +2.92%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.546"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.646.36"
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("E10").Value = "  +6.92%  "
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.86"
$ws.Range("E14").Value = "  +2.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000188"
$ws.Range("E15").Value = "  +3.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.127.48"
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.326.53"
$ws.Range("E17").Value = "  +2.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.656.03"
$ws.Range("E18").Value = "  +1.86%  "
$ws.Range("E19").Value = "  +3.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "366.24"
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("E23").Value = "  +2.60%  "
$ws.Range("E24").Value = "  +2.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.35"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.01"
$ws.Range("E27").Value = "  +0.97%  "
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("E29").Value = "  +5.87%  "
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "574.39"
$ws.Range("E31").Value = "  -1.50%  "
$ws.Range("E32").Value = "  +5.00%  "
$ws.Range("E33").Value = "  +4.56%  "
$ws.Range("E34").Value = "  +2.64%  "
$ws.Range("E35").Value = "  +3.11%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  +3.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "159.55"
$ws.Range("E38").Value = "  +1.89%  "
$ws.Range("E39").Value = "  +4.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.22"
$ws.Range("E40").Value = "  +1.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.368"
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.39"
$ws.Range("E42").Value = "  +3.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.75"
$ws.Range("E43").Value = "  +3.83%  "
$ws.Range("E44").Value = "  +2.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0320"
$ws.Range("E45").Value = "  +12.23%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.55"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "157.11"
$ws.Range("E48").Value = "  +2.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.74"
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.71"
$ws.Range("E50").Value = "  +1.67%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.87"
$ws.Range("E51").Value = "  +2.55%  "
